# Update report table: insert a new data row (IL / Full / December 2017 /
# "Order of columns differs" / Aku 06/21/2022) above the existing 2019 row,
# and keep the autofilter / named ranges / selection in sync with the
# now-larger table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 5 (existing rows 5-8 shift down to 6-9), inheriting the
# formatting of the row above, exactly like a manual Excel "Insert Row".
$ws.Rows("5:5").Insert()
$ws.Rows("5:5").RowHeight = 14.25

# The new row was pasted in from another sheet/source, so it carries an
# explicit black font color (rather than the workbook's theme color) on
# every cell from A to K.
$rng = $ws.Range("A5:K5")
$rng.Font.Color = 0

# Fill in the new row's data.
$ws.Range("A5").Value = 2018
$ws.Range("B5").Value = "IL"
$ws.Range("C5").Value = "Full"
$ws.Range("D5").Value = "December"
$ws.Range("E5").Value = 2017
$ws.Range("F5").Value = "Order of columns differs"
$ws.Range("H5").Value = "Aku 06/21/2022"

# Grow the worksheet autofilter to cover the new row.
$ws.AutoFilterMode = $false
$ws.Range("A1:F9").AutoFilter()

# Keep the hidden filter-tracking named ranges (used by the custom filter
# view) in sync with the bigger table too.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$F`$9"
  }
  if ($n.Name -eq "Sheet1!Z_00CB0BAB_A7CD_4522_9599_5EF0CAC2BF84_.wvu.FilterData") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$K`$9"
  }
}

# Restore the cursor to where the editor left it.
$ws.Range("I6").Select()
